$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 68, pushing the existing rows 68-69 down to 69-70.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly price record.
$ws.Range("A68").Value = 7
$ws.Range("B68").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C68").Value = "Ñuble"
$ws.Range("D68").Value = 45041
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112001
$ws.Range("G68").Value = "Berenjena"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 50
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = 10000
$ws.Range("N68").Value = "`$/caja 60 unidades"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 167
$ws.Range("Q68").Value = 60
$ws.Range("R68").Value = "Hortaliza"
